$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Table 1 (Source Type: Statistical Institution)
# Row 11 - Enterprises density (per 1000 people)
Set-TextValue "B11" "56.13"
Set-TextValue "C11" "3.26"
Set-TextValue "D11" "59.39"

# Row 12 - Employment (% of total)
Set-TextValue "C12" "36.69"
Set-TextValue "D12" "64.19"

# Table 2 (Source Type: SME Associations)
# Row 33 - Enterprises density (per 1000 people)
Set-TextValue "B33" "38.48"
Set-TextValue "C33" "3.33"
Set-TextValue "D33" "41.81"

# Row 34 - Employment (% of total)
Set-TextValue "B34" "24.89"
Set-TextValue "D34" "62.99"

# Row 36 - Enterprises (% of total)
Set-TextValue "B36" "91.79"
Set-TextValue "C36" "7.94"
Set-TextValue "D36" "99.73"

# Row 40 - Value added to the economy (% of total)
Set-TextValue "C40" "37.48"
Set-TextValue "D40" "59.38"
